$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PoiFormulaHelperTest")
$ws2 = $wb.Worksheets.Item("Data")

# Add new row 9 to the PoiFormulaHelperTest sheet
$ws1.Range("A9").Formula = "=SUM(Data!A1:D5)+SUM(Data!A1:D5)"
$ws1.Range("B9").Value = "Multiple Function Eval"

# Update selections to match the post-edit state
$ws2.Range("A1:D5").Select()
$ws1.Activate()
$ws1.Range("A9").Select()
